$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the bf2592bb row
# now has an actual handback timestamp instead of the placeholder status text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 08:59:54"

# zh-cn sheet: fill in the Correspond Handoff/Handback datetimes for the
# bf2592bb row (row 3) now that the handback has been generated.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-31 08:59:49"
$wsZhCn.Range("K3").Value = "2016-08-31 09:00:51"

# de-de sheet: same, for its own timestamps.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-31 08:59:54"
$wsDeDe.Range("K3").Value = "2016-08-31 09:01:07"
